$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G4").Value = 1.36
$ws.Range("G5").Value = 2.38
$ws.Range("M29").Value = 1.02
$ws.Range("O29").Value = 1.17
$ws.Range("M30").Value = 1.03
$ws.Range("O30").Value = 1.22
$ws.Range("G35").Value = 1.57
$ws.Range("H35").Value = 3.8
$ws.Range("I35").Value = 6
$ws.Range("J35").Value = 2.2
$ws.Range("M35").Value = 1.06
$ws.Range("N35").Value = 10
$ws.Range("S35").Value = 1.4
$ws.Range("T35").Value = 2.75
$ws.Range("U35").Value = 2.1
$ws.Range("V35").Value = 1.67
$ws.Range("Z35").Value = 11
$ws.Range("AE35").Value = 21
$ws.Range("AK35").Value = 67
$ws.Range("AL35").Value = 51
$ws.Range("AN35").Value = 3.4
$ws.Range("AO35").Value = 8
$ws.Range("AQ35").Value = 26
$ws.Range("AT35").Value = 2.75
$ws.Range("AX35").Value = 34
$ws.Range("Q36").Value = 1.88
$ws.Range("R36").Value = 1.98
$ws.Range("R40").Value = 1.36
$ws.Range("Q41").Value = 1.75
$ws.Range("Q42").Value = 1.73
$ws.Range("R42").Value = 2.08
